# Update the "Förändrad" (Changed) date column (C) for every data row on the
# sheet: the stored serial date moves from 45202 (2023-10-03) to
# 45203 (2023-10-04) for rows 2 through 332.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C332").Value = 45203
